# "Add Panels" is the active/selected sheet in this workbook already.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (P405D / FIM / 3.1 / 5) is removed entirely; row 11 (MX4000 / CPU 800
# / FIM / 2.5 / 5) shifts up to become the new row 10.
$ws.Rows.Item(10).Delete()

# Leave the selection where the author left off editing.
$ws.Range("B9").Select()
